$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the resistor-divider input row (row 7) ---
# Values shifted: a new intermediate data point (E7) was inserted into the
# series and the remaining inputs were rescaled down.
$ws.Range("D7").Value = 150
$ws.Range("E7").Value = 220
$ws.Range("F7").Value = 680
$ws.Range("G7").Value = 1000
$ws.Range("H7").Value = 2200

# --- Update the ADC max value ---
$ws.Range("C10").Value = 4095

# --- Update the active window view/selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("M7").Select()
